$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.288.07'
$ws.Range("E2").Value = '  +6.49%  '
$ws.Range("D3").Value = '3.557.69'
$ws.Range("E3").Value = '  +3.80%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '419.89'
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").Value = '132.43'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("D7").Value = '0.661'
$ws.Range("E7").Value = '  +6.07%  '
$ws.Range("D8").Value = '3.544.44'
$ws.Range("E8").Value = '  +3.63%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +8.27%  '
$ws.Range("E11").Value = '  +21.22%  '
$ws.Range("D12").Value = '0.0000287'
$ws.Range("E12").Value = '  +33.22%  '
$ws.Range("D13").Value = '43.50'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("E14").Value = '  +8.69%  '
$ws.Range("D15").Value = '4.120.36'
$ws.Range("E15").Value = '  +3.74%  '
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '20.62'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '3.563.63'
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").Value = '1.12'
$ws.Range("E20").Value = '  +3.95%  '
$ws.Range("D21").Value = '66.212.94'
$ws.Range("E21").Value = '  +6.31%  '
$ws.Range("D22").Value = '450.69'
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("D23").Value = '90.54'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").Value = '3.25'
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("D25").Value = '13.22'
$ws.Range("E25").Value = '  -2.35%  '
$ws.Range("D26").Value = '3.40'
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  -4.31%  '
$ws.Range("D28").Value = '34.39'
$ws.Range("E28").Value = '  +3.77%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").Value = '4.84'
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '12.50'
$ws.Range("E30").Value = '  +4.52%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '2.79'
$ws.Range("E31").Value = '  +5.43%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  +5.59%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = '7.33'
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.162'
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '39.25'
$ws.Range("E36").Value = '  -3.51%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '57.63'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0509'
$ws.Range("E38").Value = '  +4.38%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0735'
$ws.Range("E39").Value = '  +41.86%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.149'
$ws.Range("E40").Value = '  +10.90%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '3.07'
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  +4.68%  '
$ws.Range("D44").Value = '4.48'
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '147.37'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").Value = '3.28'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = '0.313'
$ws.Range("E47").Value = '  -4.02%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '2.01'
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '2.33'
$ws.Range("E49").Value = '  -6.46%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.147'
$ws.Range("E50").Value = '  +6.17%  '
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = '15.81'
$ws.Range("E51").Value = '  -3.89%  '
